# Atualizações dados 17/07 00h
#
# data/matches.xlsx, sheet "Sheet1":
#   - dt_insertion (col P) refreshed to the new scrape timestamp for a
#     handful of still-"postponed"/"notstarted" rows (13, 23, 33, 150, 153).
#   - row 150 also got its match date (col K) pushed back a day
#     (2024-07-09 -> 2024-07-10).
#   - rodada 17 (rows 163-169) got re-ordered: the match that already
#     finished (id 12117139, São Paulo x Rio de Janeiro @ Neo Química
#     Arena) dropped to the bottom of the block and the remaining six
#     postponed/not-started fixtures each shifted up one row, all
#     stamped with the new dt_insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple dt_insertion-only refreshes ---
$ws.Range("P13").Value = 45490.00259259259
$ws.Range("P23").Value = 45490.00259259259
$ws.Range("P33").Value = 45490.00259259259
$ws.Range("P153").Value = 45490.00259259259

# --- Row 150: date pushed a day later, plus dt_insertion refresh ---
$ws.Range("K150").Value = 45483
$ws.Range("P150").Value = 45490.00259259259

# --- Rodada 17 block (rows 163-169) reshuffled -------------------------
# New row 163 = old row 164 (Rio de Janeiro x Maracanã, postponed)
$ws.Range("A163").Value = 12117137
$ws.Range("B163").Value = 1961
$ws.Range("C163").Value = 1967
$ws.Range("G163").Value = "postponed"
$ws.Range("H163").Value = ""
$ws.Range("I163").Value = ""
$ws.Range("J163").Value = ""
$ws.Range("K163").Value = 45490
$ws.Range("L163").Value = "Rio de Janeiro"
$ws.Range("M163").Value = "Estádio do Maracanã"
$ws.Range("N163").Value = ""
$ws.Range("O163").Value = ""
$ws.Range("P163").Value = 45490.00259259259

# New row 164 = old row 165 (Porto Alegre x Beira-Rio, postponed)
$ws.Range("A164").Value = 12117134
$ws.Range("B164").Value = 1966
$ws.Range("C164").Value = 5981
$ws.Range("L164").Value = "Porto Alegre"
$ws.Range("M164").Value = "Estádio Beira-Rio"
$ws.Range("P164").Value = 45490.00260416666

# New row 165 = old row 166 (Goiânia x Antônio Accioly, notstarted)
$ws.Range("A165").Value = 12117142
$ws.Range("B165").Value = 7314
$ws.Range("C165").Value = 1974
$ws.Range("G165").Value = "notstarted"
$ws.Range("L165").Value = "Goiânia"
$ws.Range("M165").Value = "Estádio Antônio Accioly"
$ws.Range("P165").Value = 45490.00260416666

# New row 166 = old row 167 (São Paulo x Cícero Pompeu de Toledo, notstarted)
$ws.Range("A166").Value = 12117140
$ws.Range("B166").Value = 1981
$ws.Range("C166").Value = 5926
$ws.Range("L166").Value = "São Paulo"
$ws.Range("M166").Value = "Estádio Cícero Pompeu de Toledo"
$ws.Range("P166").Value = 45490.00260416666

# New row 167 = old row 168 (Rio de Janeiro x Nilton Santos, notstarted)
$ws.Range("A167").Value = 12117138
$ws.Range("B167").Value = 1958
$ws.Range("C167").Value = 1963
$ws.Range("L167").Value = "Rio de Janeiro"
$ws.Range("M167").Value = "Estádio Nilton Santos"
$ws.Range("P167").Value = 45490.00260416666

# New row 168 = old row 169 (Fortaleza x Castelão, notstarted)
$ws.Range("A168").Value = 12117141
$ws.Range("B168").Value = 2020
$ws.Range("C168").Value = 1962
$ws.Range("L168").Value = "Fortaleza"
$ws.Range("M168").Value = "Estádio Castelão"
$ws.Range("P168").Value = 45490.00260416666

# New row 169 = old row 163 (the finished São Paulo x Rio de Janeiro match)
$ws.Range("A169").Value = 12117139
$ws.Range("B169").Value = 1957
$ws.Range("C169").Value = 1984
$ws.Range("G169").Value = "finished"
$ws.Range("H169").Value = 788983
$ws.Range("I169").Value = 784889
$ws.Range("J169").Value = 791416
$ws.Range("K169").Value = 45489
$ws.Range("L169").Value = "São Paulo"
$ws.Range("M169").Value = "Neo Química Arena"
$ws.Range("N169").Value = 2
$ws.Range("O169").Value = 1
$ws.Range("P169").Value = 45490.00259259259
